$d = $word.ActiveDocument

# The document's auto-tracked "last edit" bookmark currently sits mid-word in
# the first edit below; remove it now and we will drop it back in at the
# location of the very last edit once all the text changes are made.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# 1) The sentence "Revising and extending previous drafts of nuclear and
#    chloroblast genome sequences of Q. lobata " had accidentally been split
#    mid-word ("...extending p" / "revious drafts..."). Re-type it as one
#    continuous run.
$r1 = $d.Content
$r1.Find.Execute("Revising and extending previous drafts of nuclear and chloroblast genome sequences of Q. lobata ", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Revising and extending previous drafts of nuclear and chloroblast genome sequences of Q. lobata ", 2) | Out-Null

# 2) Update the closing sentence of the document: wrap the
#    "possibly identify ... climate change" clause in brackets and append a
#    parenthetical aside before the final period.
$r2 = $d.Content
$r2.Find.Execute(" and possibly identify new gene sequences responsible for the plant" + [char]0x2019 + "s adaptation  to climate change.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    " and [possibly identify new gene sequences responsible for the plant" + [char]0x2019 + "s adaptation  to climate change] (actually the part before the brackets, but technically we might).", 2) | Out-Null

# 3) Word keeps the "_GoBack" bookmark pinned to the location of the most
#    recent edit -- that is now right before the final period we just typed.
$r3 = $d.Content
$r3.Find.Execute("might).") | Out-Null
$goBackPos = $r3.End - 1
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
